{"js": "// Office.js (Word JavaScript API) script.\n// Applies the two textual edits described by the diff:\n//   1. The paragraph that used to read just \"paragrafo2.\" now starts with\n//      a new \"proposicaoOuRelatorio \" placeholder run.\n//   2. The paragraph \"Constam como documentos da proposta: ...\" drops the\n//      hard-coded \"e (b) teste documentos a ser enviados\" text in favor of\n//      a \"documentosEnviados\" placeholder.\n\nconst body = context.document.body;\n\n// --- 1) Insert \"proposicaoOuRelatorio \" before \"paragrafo2\" -----------\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const para = paragraphs.items[i];\n  if (para.text.indexOf(\"paragrafo2\") === 0) {\n    para.insertText(\"proposicaoOuRelatorio \", Word.InsertLocation.start);\n    break;\n  }\n}\nawait context.sync();\n\n// --- 2) Replace the hard-coded placeholder with \"documentosEnviados\" --\nconst oldPhrase = \"e (b) teste documentos a ser enviados\";\nconst found = body.search(oldPhrase, { matchCase: true });\nfound.load(\"items\");\nawait context.sync();\n\nif (found.items.length > 0) {\n  found.items[0].insertText(\"documentosEnviados\", Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Word COM interop (PowerShell-style) script.\n# Applies the two textual edits described by the diff:\n#   1. The paragraph that used to read just \"paragrafo2.\" now starts with\n#      a new \"proposicaoOuRelatorio \" placeholder run.\n#   2. The paragraph \"Constam como documentos da proposta: ...\" drops the\n#      hard-coded \"e (b) teste documentos a ser enviados\" text in favor of\n#      a \"documentosEnviados\" placeholder.\n\n$d = $word.ActiveDocument\n\n# --- 1) Insert \"proposicaoOuRelatorio \" before \"paragrafo2\" -----------\n$paragraphs = $d.Paragraphs\nfor ($i = 1; $i -le $paragraphs.Count; $i++) {\n    $para = $paragraphs.Item($i)\n    if ($para.Range.Text.StartsWith(\"paragrafo2\")) {\n        $para.Range.InsertBefore(\"proposicaoOuRelatorio \")\n        break\n    }\n}\n\n# --- 2) Replace the hard-coded placeholder with \"documentosEnviados\" --\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"e (b) teste documentos a ser enviados\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"documentosEnviados\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n"}
